$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "Based on old employee details & keep checking ..."
#   Collapse " old employee details " + "& " + " keep" + " checking..." (which
#   were split across 4 runs around two w:proofErr gramStart/gramEnd markers)
#   into just two runs:
#     " old employee details & "
#     " keep checking the skills set and any search related to new job we
#       can capture him the employe who can chance of leaving job."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(4)
$p1Text = $p1.Range.Text
$anchor1 = "Based on"
$afterAnchor1 = $p1.Range.Start + $anchor1.Length

$splitMarker1 = "&  keep"
$splitIdx1 = $p1Text.IndexOf($splitMarker1)
$ampEnd = $p1.Range.Start + $splitIdx1 + 2   # position right after "& "

$firstRunEnd = $afterAnchor1 + (" old employee details & ").Length
$firstChunk = $d.Range($afterAnchor1, $ampEnd)
$firstChunk.Text = " old employee details & "

$paraEnd1 = $p1.Range.End - 1
$secondChunk = $d.Range($ampEnd, $paraEnd1)
$secondChunk.Text = " keep checking the skills set and any search related to new job we can capture him the employe who can chance of leaving job."

# ---------------------------------------------------------------------------
# Change 2: paragraph "Company has old employees data who those are resigned
#   employees, based on that we can predict future. So Machine Learning is
#   best choice."
#   Collapse " old " + "employees" + " data who ... future. " + "So" + " "
#   (split around two w:proofErr gramStart/gramEnd markers) into a single run:
#     " old employees data who those are resigned employees, based on that
#       we can predict future. So "
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(7)
$p2Text = $p2.Range.Text
$idxOld = $p2Text.IndexOf(" old ")
$idxML = $p2Text.IndexOf("Machine Learning")

$start2 = $p2.Range.Start + $idxOld
$end2 = $p2.Range.Start + $idxML
$midChunk = $d.Range($start2, $end2)
$midChunk.Text = " old employees data who those are resigned employees, based on that we can predict future. So "

# ---------------------------------------------------------------------------
# Change 3: "IdentifyChangeOfEmployeeResign" -> split into three runs:
#   "IdentifyCh" + "anc" + "eOfEmployeeResign"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(15)
$p3Text = $p3.Range.Text
$word3 = "IdentifyChangeOfEmployeeResign"
$idx3 = $p3Text.IndexOf($word3)
$start3 = $p3.Range.Start + $idx3

$seg1End = $start3 + 10   # "IdentifyCh"
$seg2End = $seg1End + 3   # "anc"
$seg3End = $seg2End + 17  # "eOfEmployeeResign"

$r3c = $d.Range($seg2End, $seg3End)
$r3c.Text = "eOfEmployeeResign"
$r3b = $d.Range($seg1End, $seg2End)
$r3b.Text = "anc"
$r3a = $d.Range($start3, $seg1End)
$r3a.Text = "IdentifyCh"

# Force the middle segment to serialize as an independent <w:r> (Word's COM
# object model has no direct "split run" primitive; toggling a character
# property is the only way to make the engine keep adjoining same-format
# text as separate runs instead of re-coalescing it into one).
$mid = $d.Range($seg1End, $seg2End)
$mid.Bold = 1
$mid.Bold = 0
